$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.069.56"
$ws.Range("E2").Value = "  +3.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.804.97"
$ws.Range("E3").Value = "  +4.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.01"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9983"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5549"
$ws.Range("E7").Value = "  +14.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3810"
$ws.Range("E8").Value = "  +8.79%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07624"
$ws.Range("E9").Value = "  +4.94%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.14"
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.136"
$ws.Range("E11").Value = "  +8.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9978"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.20"
$ws.Range("E13").Value = "  +6.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.217"
$ws.Range("E14").Value = "  +5.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.801.84"
$ws.Range("E15").Value = "  +4.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.161"
$ws.Range("E16").Value = "  +3.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.95"
$ws.Range("E17").Value = "  +5.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001084"
$ws.Range("E18").Value = "  +4.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06496"
$ws.Range("E19").Value = "  +1.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9981"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.20"
$ws.Range("E21").Value = "  +3.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.971"
$ws.Range("E22").Value = "  +4.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.113.34"
$ws.Range("E23").Value = "  +3.69%  "
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.104"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.65"
$ws.Range("E26").Value = "  +3.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.18"
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.379"
$ws.Range("E28").Value = "  +14.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.006.53"
$ws.Range("E29").Value = "  +4.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.88"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.150"
$ws.Range("E31").Value = "  +9.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1038"
$ws.Range("E32").Value = "  +10.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.737"
$ws.Range("E33").Value = "  +6.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.598"
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02310"
$ws.Range("E35").Value = "  +5.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2127"
$ws.Range("E36").Value = "  +6.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.675"
$ws.Range("E37").Value = "  +15.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.54"
$ws.Range("E38").Value = "  +5.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.026"
$ws.Range("E39").Value = "  +5.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06040"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6312"
$ws.Range("E41").Value = "  +5.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9977"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("B43").Value = "WEMIXTOKEN"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.399"
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.152"
$ws.Range("E44").Value = "  +4.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.44"
$ws.Range("E45").Value = "  +4.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5921"
$ws.Range("E46").Value = "  +4.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.675"
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.05"
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("E49").Value = "  +4.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.137"
$ws.Range("E50").Value = "  +3.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06790"
$ws.Range("E51").Value = "  +2.13%  "
